$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.236.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.675.35'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.37'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5279'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2657'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06286'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.35'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.01%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.669.05'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.467'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5628'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008033'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.045.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.827'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '187.98'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.43'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.220'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.84'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1259'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.597'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.97'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06251'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.365'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.286'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.517'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.437'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.638'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.004'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6071'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.410'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.748'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.118'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01622'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.102.82'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8717'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.99'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.827.83'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000111'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.23'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.007'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.031'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05236'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.984'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.52%  '
